$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = "2025-10-12 12:41:48"
}
